$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.907.76"
$ws.Range("E2").Value = "  +5.58%  "
$ws.Range("D3").Value = "3.647.75"
$ws.Range("E3").Value = "  +5.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").Value = "3.641.95"
$ws.Range("E8").Value = "  +5.24%  "
$ws.Range("E10").Value = "  +6.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.677"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("E13").Value = "  +6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.75%  "
$ws.Range("D15").Value = "4.234.78"
$ws.Range("E15").Value = "  +5.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.04%  "
$ws.Range("D17").Value = "3.653.01"
$ws.Range("E17").Value = "  +5.48%  "
$ws.Range("D18").Value = "70.927.83"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("E19").Value = "  +4.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.121"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +5.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +5.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("E32").Value = "  +9.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "618.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.15%  "
$ws.Range("D37").Value = "0.0₃0833"
$ws.Range("E37").Value = "  +9.58%  "
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "3.331.37"
$ws.Range("E42").Value = "  +3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.49%  "
$ws.Range("E44").Value = "  +6.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0457"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.84%  "
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.43%  "
